# Update "想去人数" (F column) counts across the four sheets to reflect
# the newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 768
$ws.Range("F3").Value = 2783
$ws.Range("F4").Value = 1329
$ws.Range("F7").Value = 586
$ws.Range("F9").Value = 602
$ws.Range("F11").Value = 85
$ws.Range("F12").Value = 11582
$ws.Range("F13").Value = 6598
$ws.Range("F21").Value = 74
$ws.Range("F23").Value = 924
$ws.Range("F24").Value = 3641
$ws.Range("F25").Value = 56
$ws.Range("F28").Value = 166
$ws.Range("F29").Value = 313
$ws.Range("F31").Value = 265
$ws.Range("F32").Value = 297
$ws.Range("F33").Value = 5004
$ws.Range("F35").Value = 1233
$ws.Range("F36").Value = 229
$ws.Range("F37").Value = 432
$ws.Range("F38").Value = 194
$ws.Range("F39").Value = 534

# ---- 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 23
$ws.Range("F11").Value = 3671

# ---- 本地生活 (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9033
$ws.Range("F4").Value = 1822

# ---- 全部类型 (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9033
$ws.Range("F4").Value = 1822
$ws.Range("F5").Value = 768
$ws.Range("F6").Value = 2783
$ws.Range("F9").Value = 23
$ws.Range("F10").Value = 1329
$ws.Range("F13").Value = 586
$ws.Range("F16").Value = 602
$ws.Range("F18").Value = 85
$ws.Range("F19").Value = 11582
$ws.Range("F20").Value = 3671
$ws.Range("F29").Value = 74
$ws.Range("F31").Value = 924
$ws.Range("F32").Value = 3641
$ws.Range("F33").Value = 56
$ws.Range("F35").Value = 166
$ws.Range("F36").Value = 313
$ws.Range("F37").Value = 265
$ws.Range("F41").Value = 1233
$ws.Range("F42").Value = 229
$ws.Range("F43").Value = 194
$ws.Range("F44").Value = 534

$wb.Save()
